# edit.ps1 - updates the cryptos list worksheet to the figures captured by the
# "Wed Jun  7 21:41:33 UTC 2023" GitHub Actions refresh run.
#
# All the data cells in this sheet are plain text (t="inlineStr" in the
# original OOXML) even though many of their contents look numeric (e.g.
# "1.001", "16.90", "0.05880"). If we just assign a numeric-looking string to
# Range.Value, Excel happily reinterprets it as a real number and we lose the
# original text formatting (trailing zeros, thousands-dot grouping, etc).
# To keep every updated cell as genuine text - and to avoid leaving a stray
# NumberFormat behind on the cell - Set-CellText temporarily flips the
# cell's number format to Text ("@"), assigns the literal string, then resets
# the cell's Style back to "Normal" so the cell ends up with default
# formatting and a plain text value, just like the original cells.
#
# NOTE: this runtime's function-call binding only works reliably with
# positional arguments (named `-param value` args came back empty in
# testing), so Set-CellText is always called positionally below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($wsArg, $addr, $value)

    $range = $wsArg.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-CellText $ws 'D2' '26.440.27'
Set-CellText $ws 'E2' '  -2.07%  '
Set-CellText $ws 'D3' '1.839.60'
Set-CellText $ws 'E3' '  -2.02%  '
Set-CellText $ws 'D4' '1.001'
Set-CellText $ws 'E4' '  +0.09%  '
Set-CellText $ws 'D5' '260.28'
Set-CellText $ws 'E5' '  -7.50%  '
Set-CellText $ws 'D6' '1.001'
Set-CellText $ws 'E6' '  +0.01%  '
Set-CellText $ws 'D7' '0.5218'
Set-CellText $ws 'E7' '  -0.92%  '
Set-CellText $ws 'D8' '0.3227'
Set-CellText $ws 'E8' '  -8.60%  '
Set-CellText $ws 'D9' '0.06737'
Set-CellText $ws 'E9' '  -4.24%  '
Set-CellText $ws 'D10' '18.87'
Set-CellText $ws 'E10' '  -7.15%  '
Set-CellText $ws 'D11' '0.7685'
Set-CellText $ws 'E11' '  -5.66%  '
Set-CellText $ws 'D12' '0.07681'
Set-CellText $ws 'E12' '  -1.32%  '
Set-CellText $ws 'D13' '1.831.01'
Set-CellText $ws 'E13' '  -2.41%  '
Set-CellText $ws 'D14' '89.16'
Set-CellText $ws 'E14' '  -1.36%  '
Set-CellText $ws 'D15' '5.026'
Set-CellText $ws 'E15' '  -3.54%  '
Set-CellText $ws 'D16' '1.001'
Set-CellText $ws 'E16' '  +0.08%  '
Set-CellText $ws 'D17' '14.11'
Set-CellText $ws 'E17' '  -3.25%  '
Set-CellText $ws 'E18' '  +0.03%  '
Set-CellText $ws 'D19' '0.000007859'
Set-CellText $ws 'E19' '  -3.91%  '
Set-CellText $ws 'D20' '26.483.03'
Set-CellText $ws 'E20' '  -1.98%  '
Set-CellText $ws 'D21' '2.087.16'
Set-CellText $ws 'E21' '  -0.90%  '
Set-CellText $ws 'D22' '4.539'
Set-CellText $ws 'E22' '  -4.64%  '
Set-CellText $ws 'D23' '9.461'
Set-CellText $ws 'E23' '  -7.08%  '
Set-CellText $ws 'D24' '5.913'
Set-CellText $ws 'E24' '  -5.05%  '
Set-CellText $ws 'D25' '2.335'
Set-CellText $ws 'E25' '  -2.02%  '
Set-CellText $ws 'D26' '144.32'
Set-CellText $ws 'E26' '  -1.29%  '
Set-CellText $ws 'D27' '1.648'
Set-CellText $ws 'E27' '  -1.60%  '
Set-CellText $ws 'D28' '16.90'
Set-CellText $ws 'E28' '  -3.73%  '
Set-CellText $ws 'D29' '111.17'
Set-CellText $ws 'E29' '  -1.85%  '
Set-CellText $ws 'D30' '4.179'
Set-CellText $ws 'E30' '  -4.55%  '
Set-CellText $ws 'B31' 'Filecoin'
Set-CellText $ws 'C31' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws 'D31' '4.125'
Set-CellText $ws 'E31' '  -5.65%  '
Set-CellText $ws 'B32' 'Stellar'
Set-CellText $ws 'C32' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-CellText $ws 'D32' '0.08759'
Set-CellText $ws 'E32' '  -1.42%  '
Set-CellText $ws 'D33' '0.04837'
Set-CellText $ws 'E33' '  -1.23%  '
Set-CellText $ws 'D34' '1.133'
Set-CellText $ws 'E34' '  -3.48%  '
Set-CellText $ws 'D35' '2.851'
Set-CellText $ws 'E35' '  -0.90%  '
Set-CellText $ws 'D36' '0.6841'
Set-CellText $ws 'E36' '  -7.87%  '
Set-CellText $ws 'D37' '3.109'
Set-CellText $ws 'E37' '  -5.60%  '
Set-CellText $ws 'D38' '0.01787'
Set-CellText $ws 'E38' '  -5.17%  '
Set-CellText $ws 'D39' '2.219'
Set-CellText $ws 'E39' '  -7.70%  '
Set-CellText $ws 'D40' '0.4924'
Set-CellText $ws 'E40' '  -6.99%  '
Set-CellText $ws 'D41' '112.47'
Set-CellText $ws 'E41' '  -3.96%  '
Set-CellText $ws 'D42' '0.8955'
Set-CellText $ws 'E42' '  -8.69%  '
Set-CellText $ws 'D43' '6.174'
Set-CellText $ws 'D44' '1.000'
Set-CellText $ws 'E44' '  +0.04%  '
Set-CellText $ws 'D45' '7.742'
Set-CellText $ws 'E45' '  -5.47%  '
Set-CellText $ws 'D46' '0.4199'
Set-CellText $ws 'E46' '  -8.43%  '
Set-CellText $ws 'E47' '  -7.77%  '
Set-CellText $ws 'D48' '9.072'
Set-CellText $ws 'E48' '  -4.07%  '
Set-CellText $ws 'B49' 'Cronos'
Set-CellText $ws 'C49' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText $ws 'D49' '0.05880'
Set-CellText $ws 'E49' '  -1.10%  '
Set-CellText $ws 'B50' 'Elrond'
Set-CellText $ws 'C50' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-CellText $ws 'D50' '35.41'
Set-CellText $ws 'E50' '  -3.62%  '
Set-CellText $ws 'B51' 'NEARProtocol'
Set-CellText $ws 'C51' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-CellText $ws 'D51' '1.418'
Set-CellText $ws 'E51' '  -6.69%  '
